$wb = $excel.ActiveWorkbook

# --- Global sheet ---
$wsGlobal = $wb.Worksheets.Item("Global")
$wsGlobal.Range("B3").Value = "SIS3316_125"
$wsGlobal.Range("D3").Value = 4
$wsGlobal.Range("G3").Value = "LEMO"
$wsGlobal.Range("H3").Value = 25

# --- Rename sheets ---
$wsHard = $wb.Worksheets.Item("Hard_Sis3316_125_0xCC")
$wsHard.Name = "Hard_SIS3316_125"

$wsCal = $wb.Worksheets.Item("Cal_Sis3316_125_0xCC")
$wsCal.Name = "Cal_SIS3316_125"

# --- Hard_SIS3316_125 sheet (channel hardware parameters) ---
# Row 4
$wsHard.Range("E4").Value = 10
$wsHard.Range("G4").Value = 700
$wsHard.Range("H4").Value = 220
$wsHard.Range("J4").Value = 100
$wsHard.Range("K4").Value = 200
$wsHard.Range("L4").Value = 260000
$wsHard.Range("P4").Value = 800
$wsHard.Range("AA4").Value = 0
$wsHard.Range("AB4").Value = 55000
$wsHard.Range("AD4").Value = 1
$wsHard.Range("AI4").Value = 800000

# Row 5
$wsHard.Range("E5").Value = 10
$wsHard.Range("G5").Value = 620
$wsHard.Range("J5").Value = 100
$wsHard.Range("K5").Value = 200
$wsHard.Range("L5").Value = 500000
$wsHard.Range("AA5").Value = 0
$wsHard.Range("AB5").Value = 55400
$wsHard.Range("AD5").Value = 1
$wsHard.Range("AI5").Value = 800000

# Row 6
$wsHard.Range("E6").Value = 10
$wsHard.Range("G6").Value = 640
$wsHard.Range("J6").Value = 100
$wsHard.Range("K6").Value = 200
$wsHard.Range("L6").Value = 260000
$wsHard.Range("AA6").Value = 0
$wsHard.Range("AB6").Value = 55400
$wsHard.Range("AD6").Value = 1
$wsHard.Range("AI6").Value = 800000

# Row 7
$wsHard.Range("E7").Value = 10
$wsHard.Range("G7").Value = 640
$wsHard.Range("J7").Value = 100
$wsHard.Range("K7").Value = 200
$wsHard.Range("L7").Value = 500000
$wsHard.Range("AA7").Value = 0
$wsHard.Range("AB7").Value = 55500
$wsHard.Range("AD7").Value = 1
$wsHard.Range("AI7").Value = 800000

# --- Cal_SIS3316_125 sheet (calibration factors) ---
$wsCal.Range("I3").Value = 0.0472684873949579
$wsCal.Range("J3").Value = -39.32352941176467

$wsCal.Range("I4").Value = 0.0438307792207792
$wsCal.Range("J4").Value = -38.2000000000001

$wsCal.Range("I5").Value = 0.0241069285714285
$wsCal.Range("J5").Value = -46.79500000000002

$wsCal.Range("I6").Value = 0.0251113839285714
$wsCal.Range("J6").Value = -29.84375000000003

# --- Condition sheet: divide column K detection limits by 40 ---
$wsCond = $wb.Worksheets.Item("Condition")

$wsCond.Range("K3").Value = 0.000000004
$wsCond.Range("K6").Value = 0.00000000006666666666666667
$wsCond.Range("K9").Value = 0.000000000004444444444444445
$wsCond.Range("K12").Value = 0.000000000002222222222222222
$wsCond.Range("K15").Value = 0.000000000001111111111111111

$wsCond.Range("K18").Value = 0.000000004
$wsCond.Range("K21").Value = 0.00000000006666666666666667
$wsCond.Range("K24").Value = 0.000000000004444444444444445
$wsCond.Range("K27").Value = 0.000000000002222222222222222
$wsCond.Range("K30").Value = 0.000000000001111111111111111

$wsCond.Range("K33").Value = 0.000000004
$wsCond.Range("K36").Value = 0.00000000006666666666666667
$wsCond.Range("K39").Value = 0.000000000004444444444444445
$wsCond.Range("K42").Value = 0.000000000002222222222222222
$wsCond.Range("K45").Value = 0.000000000001111111111111111

$wsCond.Range("K48").Value = 0.000000004
$wsCond.Range("K51").Value = 0.00000000006666666666666667
$wsCond.Range("K54").Value = 0.000000000004444444444444445
$wsCond.Range("K57").Value = 0.000000000002222222222222222
$wsCond.Range("K60").Value = 0.000000000001111111111111111
